$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last edited" date cell (C1) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- MCF sheet: update capacity factor values ---
$mcf = $wb.Worksheets.Item("MCF")

$mcf.Range("B2").Value = 1
$mcf.Range("B3").Value = 1
$mcf.Range("B4").Value = 1
$mcf.Range("B6").Value = 1
$mcf.Range("B10").Value = 1
$mcf.Range("B11").Value = 1
$mcf.Range("B12").Value = 1
$mcf.Range("B13").Value = 1
$mcf.Range("B14").Value = 1
$mcf.Range("B16").Value = 1
$mcf.Range("B17").Value = 1
$mcf.Range("B18").Value = 1

# Recalculate so dependent formula cells (B19, B20, B21, B22, B24, B25) pick up new values
$excel.Calculate()

# Update the selected cell on the MCF sheet to match the diff
$mcf.Activate()
$mcf.Range("B17").Select()
